$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.895.81"
$ws.Range("E2").Value = "  +1.06%  "

$ws.Range("D3").Value = "'3.111.13"
$ws.Range("E3").Value = "  +1.27%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'576.26"
$ws.Range("E5").Value = "  -0.38%  "

$ws.Range("D6").Value = "'173.07"
$ws.Range("E6").Value = "  +3.54%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "'3.107.66"
$ws.Range("E8").Value = "  +1.23%  "

$ws.Range("E9").Value = "  -0.22%  "

$ws.Range("D10").Value = "'6.46"
$ws.Range("E10").Value = "  -2.77%  "

$ws.Range("D11").Value = "'0.153"
$ws.Range("E11").Value = "  +0.02%  "

$ws.Range("D12").Value = "'0.478"
$ws.Range("E12").Value = "  -0.32%  "

$ws.Range("E13").Value = "  -1.33%  "

$ws.Range("D14").Value = "'37.08"
$ws.Range("E14").Value = "  +0.96%  "

$ws.Range("E15").Value = "  -0.98%  "

$ws.Range("D16").Value = "'3.625.88"
$ws.Range("E16").Value = "  +1.37%  "

$ws.Range("D17").Value = "'66.906.03"
$ws.Range("E17").Value = "  +1.08%  "

$ws.Range("D18").Value = "'7.10"
$ws.Range("E18").Value = "  -1.09%  "

$ws.Range("D19").Value = "'3.111.04"
$ws.Range("E19").Value = "  +1.38%  "

$ws.Range("D20").Value = "'16.29"
$ws.Range("E20").Value = "  +1.09%  "

$ws.Range("D21").Value = "'476.41"
$ws.Range("E21").Value = "  +2.66%  "

$ws.Range("D22").Value = "'0.711"
$ws.Range("E22").Value = "  -0.02%  "

$ws.Range("D23").Value = "'7.82"
$ws.Range("E23").Value = "  +5.21%  "

$ws.Range("D24").Value = "'13.39"
$ws.Range("E24").Value = "  +4.62%  "

$ws.Range("D25").Value = "'83.74"
$ws.Range("E25").Value = "  +0.83%  "

$ws.Range("E26").Value = "  +0.24%  "

$ws.Range("D27").Value = "'1.00"

$ws.Range("E28").Value = "  -1.52%  "

$ws.Range("E29").Value = "  +0.21%  "

$ws.Range("D30").Value = "'7.90"
$ws.Range("E30").Value = "  -2.36%  "

$ws.Range("E31").Value = "  -0.28%  "

$ws.Range("D32").Value = "'28.65"
$ws.Range("E32").Value = "  +1.45%  "

$ws.Range("E33").Value = "  -1.40%  "

$ws.Range("D34").Value = "'0.0₃0941"
$ws.Range("E34").Value = "  -8.06%  "

$ws.Range("E35").Value = "  +0.05%  "

$ws.Range("E36").Value = "  -0.29%  "

$ws.Range("D37").Value = "'0.978"
$ws.Range("E37").Value = "  -1.78%  "

$ws.Range("D38").Value = "'47.49"
$ws.Range("E38").Value = "  -2.72%  "

$ws.Range("E39").Value = "  +2.59%  "

$ws.Range("D40").Value = "'49.92"
$ws.Range("E40").Value = "  -0.09%  "

$ws.Range("E41").Value = "  -1.07%  "

$ws.Range("E42").Value = "  -0.20%  "

$ws.Range("E43").Value = "  -0.93%  "

$ws.Range("D44").Value = "'2.789.33"
$ws.Range("E44").Value = "  +1.11%  "

$ws.Range("E45").Value = "  -1.67%  "

$ws.Range("D46").Value = "'377.54"
$ws.Range("E46").Value = "  -0.87%  "

$ws.Range("D47").Value = "'2.54"
$ws.Range("E47").Value = "  -12.26%  "

$ws.Range("D48").Value = "'136.06"
$ws.Range("E48").Value = "  +1.05%  "

$ws.Range("D50").Value = "'24.74"
$ws.Range("E50").Value = "  +1.20%  "

$ws.Range("E51").Value = "  -0.96%  "
